# ---------------------------------------------------------------------------
# UPS.xlsx: rename the sheet tab and populate 3 additional UPS inventory rows
# (rows 4-6), plus backfill a couple of newly-inserted columns on the existing
# rows 2-3 so every row lines up against the new shared-string table.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "UPS" -> "Sheet1"
$ws.Name = "Sheet1"

# xlPasteValues - used below to drop date-shaped text (e.g. "2024-08-05") into a
# cell without Excel reinterpreting it as a date serial number: we first land the
# text via a `="..."` formula (formula results are never auto-converted), then
# copy/paste-special-values over itself to collapse it to a plain literal value.
$xlPasteValues = -4163

function Set-TextValue($range, [string]$text) {
    $range.Formula = '="' + ($text -replace '"', '""') + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial($xlPasteValues) | Out-Null
}

# Row 2
$ws.Range("A2").Value = "UPS-001"
$ws.Range("C2").Value = "APC Smart-UPS SC 1500VA"
$ws.Range("D2").Value = "APC"
$ws.Range("E2").Value = 1500
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = "GAB-O-P3"
$ws.Range("H2").Value = "Edificio O - 3er Piso - Sala Técnica"
$ws.Range("I2").Value = "Andrés Bello"
$ws.Range("J2").Value = "Operativo"
Set-TextValue $ws.Range("K2") "2023-03-15"
$ws.Range("L2").Value = "11 cámaras Edificio O + 1 cámara PTZ"
Set-TextValue $ws.Range("M2") "2024-10-13"
$ws.Range("N2").Value = 45000
$ws.Range("O2").Value = "Cambio de batería realizado el 13/10/2024. Sistema funcionó con batería restante durante cambio."

# Row 3
$ws.Range("A3").Value = "UPS-002"
$ws.Range("C3").Value = "Tripp Lite 1000VA"
$ws.Range("D3").Value = "Tripp Lite"
$ws.Range("E3").Value = 1000
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "GAB-CFT-1"
$ws.Range("H3").Value = "CFT Prat - Sala Servidores"
$ws.Range("I3").Value = "Andrés Bello"
$ws.Range("J3").Value = "Operativo"
Set-TextValue $ws.Range("K3") "2023-05-20"
$ws.Range("L3").Value = "13 cámaras CFT Prat"
$ws.Range("O3").Value = "Programar revisión de baterías"

# Row 4
$ws.Range("A4").Value = "UPS-002"
$ws.Range("B4").Value = "UBI-006"
$ws.Range("C4").Value = "Smart-UPS 1000"
$ws.Range("D4").Value = "APC"
$ws.Range("E4").Value = 1000
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "GAB-004"
$ws.Range("H4").Value = "Campus Pucón - Recepción"
$ws.Range("I4").Value = "Campus Pucón"
$ws.Range("J4").Value = "Activo"
Set-TextValue $ws.Range("K4") "2024-08-05"
$ws.Range("L4").Value = "1 Switch + 5 cámaras"
Set-TextValue $ws.Range("M4") "2025-08-05"
$ws.Range("N4").Value = 25000
$ws.Range("O4").Value = "Protección básica"

# Row 5
$ws.Range("A5").Value = "UPS-003"
$ws.Range("B5").Value = "UBI-007"
$ws.Range("C5").Value = "Smart-UPS 2200"
$ws.Range("D5").Value = "APC"
$ws.Range("E5").Value = 2200
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = "GAB-005"
$ws.Range("H5").Value = "CFT Prat - Sala servidores"
$ws.Range("I5").Value = "Campus Angol"
$ws.Range("J5").Value = "Activo"
Set-TextValue $ws.Range("K5") "2024-05-12"
$ws.Range("L5").Value = "2 Switches + 1 NVR + 13 cámaras"
Set-TextValue $ws.Range("M5") "2025-05-12"
$ws.Range("N5").Value = 65000
$ws.Range("O5").Value = "Punto crítico - UPS redundante"

# Row 6
$ws.Range("A6").Value = "UPS-004"
$ws.Range("B6").Value = "UBI-004"
$ws.Range("C6").Value = "Back-UPS 700"
$ws.Range("D6").Value = "Tripp Lite"
$ws.Range("E6").Value = 700
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = "GAB-006"
$ws.Range("H6").Value = "Edificio L - 2do Piso"
$ws.Range("I6").Value = "Campus Principal"
$ws.Range("J6").Value = "Activo"
Set-TextValue $ws.Range("K6") "2024-06-20"
$ws.Range("L6").Value = "1 Switch + 8 cámaras"
$ws.Range("O6").Value = "Pendiente primera mantención"

